$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Row 11: Enterprises density (per 1000 people) - Source Type: Statistical Institution
Set-TextValue "B11" "68.58"
Set-TextValue "C11" "5.14"
Set-TextValue "D11" "73.72"

# Row 12: Employment (% of total) - Source Type: Statistical Institution
Set-TextValue "B12" "26.99"
Set-TextValue "C12" "34.06"
Set-TextValue "D12" "61.04"

# Row 33: Enterprises density (per 1000 people) - Source Type: SME Associations
Set-TextValue "B33" "31.53"
Set-TextValue "C33" "4.47"

# Row 34: Employment (% of total) - Source Type: SME Associations
Set-TextValue "B34" "24.42"
Set-TextValue "C34" "43.27"
Set-TextValue "D34" "67.69"

# Row 36: Enterprises (% of total) - Source Type: SME Associations
Set-TextValue "B36" "87.29"
Set-TextValue "C36" "12.37"
Set-TextValue "D36" "99.66"
